$wb = $excel.ActiveWorkbook

# --- RoboRIO sheet: rename line-follower / elevator-limit variable names ---
$ws1 = $wb.Worksheets.Item("RoboRIO")
$ws1.Range("F4").Value = "leftFollower"
$ws1.Range("F5").Value = "midFollower"
$ws1.Range("F6").Value = "rightFollower"
$ws1.Range("F7").Value = "lowLimitSwitch"

# --- PDP sheet: rename motor variable names + hatch pickup subsystem ---
$ws2 = $wb.Worksheets.Item("PDP")
$ws2.Range("J8").Value = "winchMotor"
$ws2.Range("J9").Value = "rotatorMotor"
$ws2.Range("J10").Value = "leftClimbMotor"
$ws2.Range("J11").Value = "rightClimbMotor"
$ws2.Range("J12").Value = "pickupMotor"
$ws2.Range("J13").Value = "pickupMotor"
$ws2.Range("I13").Value = "HatchPickupSubsystem"

# --- PCM sheet: rename solenoid variable names + split hatch subsystem ---
$ws4 = $wb.Worksheets.Item("PCM")
$ws4.Range("F2").Value = "pickupSol"
$ws4.Range("G2").Value = "HatchPickupSubsystem"
$ws4.Range("F3").Value = "placerSol"
$ws4.Range("G3").Value = "HatchPlacerSubsystem"
$ws4.Range("F4").Value = "detachLeftSol"
$ws4.Range("G4").Value = "HatchPlacerSubsystem"
$ws4.Range("F5").Value = "detachRightSol"
$ws4.Range("G5").Value = "HatchPlacerSubsystem"
$ws4.Range("F6").Value = "leftSol"
$ws4.Range("G6").Value = "ClimbSubsystem"
$ws4.Range("F7").Value = "rightSol"
$ws4.Range("G7").Value = "ClimbSubsystem"

# --- Re-fit the two columns whose best-fit width changed because the new text is longer/shorter ---
[void]$ws2.Columns.Item(9).AutoFit()
[void]$ws2.Columns.Item(10).AutoFit()
[void]$ws4.Columns.Item(6).AutoFit()
[void]$ws4.Columns.Item(7).AutoFit()

# --- Selection / active-sheet bookkeeping to mirror the author's final view state ---
[void]$ws1.Range("F8").Select()
[void]$ws2.Range("J10").Select()
[void]$ws4.Range("F8").Select()
[void]$ws4.Activate()
